# Archivo para hacer Facturas 3.0.xlsx - add "Domicilio / Guardar / Ubicacion a Guardar"
# columns (AX:AZ), restyle the AJ helper column, extend the autofilter / filter
# database range and add a SI/NO data-validation list on the new "Guardar" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New header cells AX1:AZ1 - reuse the existing header style (same as AM1)
#    and set their text.
# ---------------------------------------------------------------------------
$ws.Range("AM1").Copy()
$ws.Range("AX1:AZ1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AX1").Value = "Domicilio"
$ws.Range("AY1").Value = "Guardar"
$ws.Range("AZ1").Value = "Ubicación a Guardar"

# ---------------------------------------------------------------------------
# 2) New data columns AX2:AZ17 - reuse the plain body style already used by
#    the neighbouring AW column.
# ---------------------------------------------------------------------------
$ws.Range("AW2").Copy()
$ws.Range("AX2:AZ17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# AZ: save-location path only for the rows flagged "SI".
$ws.Range("AZ2:AZ10").Value = "C:\Users\ABP\Desktop\Test\"

# AY: "SI" for the first nine data rows, "NO" for the rest.
$ws.Range("AY2:AY10").Value = "SI"
$ws.Range("AY11:AY17").Value = "NO"

# ---------------------------------------------------------------------------
# 3) AJ2:AJ17 loses its extra "applyFill" formatting - restyle to match the
#    plain bordered body style (same as E2).
# ---------------------------------------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("AJ2:AJ17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Data validation (list SI / NO) on the new "Guardar" column.
# ---------------------------------------------------------------------------
$ws.Range("AY2:AY17").Validation.Delete()
$ws.Range("AY2:AY17").Validation.Add(3, 1, 1, '"SI,NO"')

# ---------------------------------------------------------------------------
# 5) Extend the AutoFilter range and the hidden _FilterDatabase name up to
#    the new AZ column.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:AZ17").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Factura!_FilterDatabase") {
        $n.RefersTo = "=Factura!`$A`$1:`$AZ`$17"
    }
}

# ---------------------------------------------------------------------------
# 6) Restore the view: default selection K17, scroll so column AK is first
#    visible, and select AY12 in the frozen bottom pane (matches the saved
#    author view).
# ---------------------------------------------------------------------------
$ws.Range("K17").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 37
$ws.Range("AY12").Select() | Out-Null
